# Update the "想去人数" (F column) figures that changed between the two
# data-collection runs. The same values were refreshed identically on the
# "展览" and "全部类型" worksheets (they hold duplicate data).

$wb = $excel.ActiveWorkbook

# Map of cell -> new value to apply on each of the affected worksheets.
$updates = @{
    "F3"  = 295
    "F4"  = 58
    "F5"  = 604
    "F6"  = 63
    "F7"  = 2080
    "F10" = 4549
    "F12" = 287
    "F15" = 136
    "F16" = 30
    "F19" = 3429
    "F20" = 85
    "F21" = 546
    "F24" = 87
    "F25" = 96
    "F28" = 67
    "F29" = 208
    "F31" = 676
    "F32" = 2076
    "F33" = 393
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cell in $updates.Keys) {
        $ws.Range($cell).Value = $updates[$cell]
    }
}
